$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A100").Value = "hello"
